$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) values for the refreshed crypto data.
# Column D values are forced to text via NumberFormat "@" so Excel does not auto-convert
# strings that look numeric (e.g. "254.00", "1.99") into actual numbers, then the style
# is reset back to Normal so no stray number-format style lingers on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.345.66"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.31"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.723"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  +10.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "254.00"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  +3.59%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.69"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  -1.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.355"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  +1.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.28"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0751"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  +5.25%  "

$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.191.74"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  +0.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.65"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  +5.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.717"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +2.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.936.62"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  +2.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.92"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  +1.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.345.92"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.57"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +4.35%  "

$ws.Range("E20").Value = "  +3.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.94"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  +1.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.06"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  +4.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.08"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  +5.48%  "

$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("E26").Value = "  +2.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.80"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.62"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  +2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.73"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  +2.35%  "

$ws.Range("E30").Value = "  +5.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.126.26"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  +19.39%  "

$ws.Range("E32").Value = "  +5.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.99"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  +14.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.64"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  +22.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0580"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  +2.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.22"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  +2.29%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.919"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  -1.85%  "

$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.47"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  +6.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0221"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  +5.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.32"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  +8.18%  "

$ws.Range("E43").Value = "  +1.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0655"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  +0.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.340.76"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.44"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  +1.99%  "

$ws.Range("E47").Value = "  +1.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.77"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  +3.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.77"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.13"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  -5.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.78"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  +6.98%  "
